$wb = $excel.ActiveWorkbook

# The nightly build stamp embedded throughout the workbook needs to move
# from the January 30 build to the February 02 build.
$oldText = "January 30 2026 16.19.47 EST"
$newText = "February 02 2026 12.49.33 EST"

# "About" sheet: version banner (A2) and recommended citation (A6)
$aboutWs = $wb.Worksheets.Item("About")

$a2 = $aboutWs.Range("A2")
$a2.Value = $a2.Value().Replace($oldText, $newText)

$a6 = $aboutWs.Range("A6")
$a6.Value = $a6.Value().Replace($oldText, $newText)

# "Boundaries and methane sources" sheet: build_version column (S) for
# every data row (2 through 16)
$dataWs = $wb.Worksheets.Item("Boundaries and methane sources")
for ($row = 2; $row -le 16; $row++) {
    $cell = $dataWs.Range("S$row")
    $cell.Value = $cell.Value().Replace($oldText, $newText)
}
